# Add BMC ports for river compute nodes (cn01-cn04) to the
# HARDWARE_MANAGEMENT and COMPUTE_NODES sheets, and update selections.

$wb = $excel.ActiveWorkbook

$wsHw = $wb.Worksheets.Item("HARDWARE_MANAGEMENT")
$wsCn = $wb.Worksheets.Item("COMPUTE_NODES")

# --- HARDWARE_MANAGEMENT: new rows 25-28 (cn01..cn04 BMC ports) ---
# Populate "Source" column first (creates cn01..cn04 shared strings),
# then the "Location" column (creates u15..u18 shared strings), matching
# the order new unique strings were appended to the workbook.
$wsHw.Range("J25").Value = "cn01"
$wsHw.Range("J26").Value = "cn02"
$wsHw.Range("J27").Value = "cn03"
$wsHw.Range("J28").Value = "cn04"

$wsHw.Range("L25").Value = "u15"
$wsHw.Range("L26").Value = "u16"
$wsHw.Range("L27").Value = "u17"
$wsHw.Range("L28").Value = "u18"

$wsHw.Range("K25:K28").Value = "x3002"
$wsHw.Range("M25:M28").Value = "bmc"
$wsHw.Range("O25:O28").Value = 1
$wsHw.Range("P25:P28").Value = "sw-leaf-bmc-001"
$wsHw.Range("Q25:Q28").Value = "x3000"
$wsHw.Range("R25:R28").Value = "u37"

$wsHw.Range("T25").Value = 11
$wsHw.Range("T26").Value = 12
$wsHw.Range("T27").Value = 13
$wsHw.Range("T28").Value = 14

# --- COMPUTE_NODES: new rows 24-27 (cn01..cn04 BMC ports) ---
$wsCn.Range("J24").Value = "cn01"
$wsCn.Range("J25").Value = "cn02"
$wsCn.Range("J26").Value = "cn03"
$wsCn.Range("J27").Value = "cn04"

$wsCn.Range("L24").Value = "u15"
$wsCn.Range("L25").Value = "u16"
$wsCn.Range("L26").Value = "u17"
$wsCn.Range("L27").Value = "u18"

$wsCn.Range("K24:K27").Value = "x3002"
$wsCn.Range("O24:O27").Value = 1
$wsCn.Range("P24:P27").Value = "sw-leaf-bmc-001"
$wsCn.Range("Q24:Q27").Value = "x3000"
$wsCn.Range("R24:R27").Value = "u37"

$wsCn.Range("T24").Value = 24
$wsCn.Range("T25").Value = 25
$wsCn.Range("T26").Value = 26
$wsCn.Range("T27").Value = 27

# --- Update selections / active sheet to match the saved UI state ---
$wsCn.Range("M24:M27").Select()

$wsHw.Activate()
$wsHw.Range("A28:XFD28").Select()
$wsHw.Range("D28").Activate()
